$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")

# Copy the formatting of the last existing data row (34) down into the
# three new rows (35-37) so styles match exactly.
$ws.Range("B34:G34").Copy()
$ws.Range("B35:G37").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new rows' content (mirrors the existing fill pattern).
for ($i = 35; $i -le 37; $i++) {
    $prev = $i - 1
    $year = 2050 + ($i - 34)

    $ws.Cells.Item($i, 2).Formula = "=Commodities!`$D`$4"
    $ws.Cells.Item($i, 3).Formula = "=C$prev"
    $ws.Cells.Item($i, 4).Formula = "=Processes!`$E`$3"
    $ws.Cells.Item($i, 5).Value = "PJ"
    $ws.Cells.Item($i, 6).Value = $year
    $ws.Cells.Item($i, 7).Formula = "=G$prev+`$G`$14*(`$I`$4)"
}

# Restore the view: scrolled down a bit with I32 selected.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I32").Select() | Out-Null
